$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values for row 32
$ws.Range("C32").Value = 70
$ws.Range("D32").Value = 80
$ws.Range("E32").Value = 80
$ws.Range("F32").Value = 90

# Update cell values for row 33
$ws.Range("C33").Value = 85
$ws.Range("D33").Value = 90
$ws.Range("E33").Value = 90
$ws.Range("F33").Value = 95

# Update cell values for row 34 (only C34 changes)
$ws.Range("C34").Value = 90

# Update cell values for row 41
$ws.Range("C41").Value = 80
$ws.Range("D41").Value = 90
$ws.Range("E41").Value = 90
$ws.Range("F41").Value = 100

# Update cell values for row 42
$ws.Range("C42").Value = 80
$ws.Range("D42").Value = 100
$ws.Range("E42").Value = 100
$ws.Range("F42").Value = 120

# Update cell values for row 43 (only C43 changes)
$ws.Range("C43").Value = 80

# Update sheet view: scroll position (topLeftCell A24) and selection (G41)
$ws.Range("G41").Select()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
